$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "shifted to selenium" — the old manually-entered lookup row (name "desai the
# god" / phone number) is no longer needed, so clear its contents but leave
# the cell formatting (style) intact.
$ws.Range("A2:B2").ClearContents()

# Select the now-empty row 3 (whole row), matching the author's last
# on-screen selection when they saved.
$null = $ws.Rows(3).Select()
